{"js": "// Replace the date heading and every three-digit-by-one-digit multiplication\n// expression in the table with their updated values, as described by the diff.\n// Each entry is [oldText, newText]; every oldText is a unique, exact match\n// within the document, so a direct search-and-replace is sufficient.\nconst replacements = [[\"2025-10-27 Monday\", \"2025-10-28 Tuesday\"], [\"852\u00d72=1704\", \"402\u00d72=804\"], [\"544\u00d73=1632\", \"129\u00d74=516\"], [\"225\u00d72=450\", \"970\u00d77=6790\"], [\"622\u00d78=4976\", \"874\u00d77=6118\"], [\"838\u00d78=6704\", \"773\u00d77=5411\"], [\"219\u00d79=1971\", \"168\u00d72=336\"], [\"448\u00d76=2688\", \"443\u00d79=3987\"], [\"364\u00d77=2548\", \"828\u00d74=3312\"], [\"399\u00d77=2793\", \"757\u00d75=3785\"], [\"904\u00d72=1808\", \"943\u00d73=2829\"], [\"669\u00d76=4014\", \"947\u00d74=3788\"], [\"128\u00d74=512\", \"465\u00d75=2325\"], [\"169\u00d77=1183\", \"959\u00d79=8631\"], [\"143\u00d74=572\", \"693\u00d73=2079\"], [\"620\u00d77=4340\", \"772\u00d76=4632\"], [\"619\u00d76=3714\", \"796\u00d78=6368\"], [\"305\u00d77=2135\", \"324\u00d73=972\"], [\"321\u00d79=2889\", \"902\u00d76=5412\"], [\"857\u00d79=7713\", \"871\u00d75=4355\"], [\"516\u00d72=1032\", \"188\u00d78=1504\"], [\"522\u00d77=3654\", \"551\u00d75=2755\"], [\"265\u00d78=2120\", \"560\u00d76=3360\"], [\"182\u00d76=1092\", \"376\u00d74=1504\"], [\"182\u00d72=364\", \"545\u00d75=2725\"], [\"878\u00d79=7902\", \"354\u00d77=2478\"]];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date heading and every three-digit-by-one-digit multiplication\n# expression in the table to its new value, per the diff. Each old string is\n# unique within the document, so Find/Replace (wdReplaceAll) is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-10-27 Monday\", \"2025-10-28 Tuesday\"),\n    @(\"852\u00d72=1704\", \"402\u00d72=804\"),\n    @(\"544\u00d73=1632\", \"129\u00d74=516\"),\n    @(\"225\u00d72=450\", \"970\u00d77=6790\"),\n    @(\"622\u00d78=4976\", \"874\u00d77=6118\"),\n    @(\"838\u00d78=6704\", \"773\u00d77=5411\"),\n    @(\"219\u00d79=1971\", \"168\u00d72=336\"),\n    @(\"448\u00d76=2688\", \"443\u00d79=3987\"),\n    @(\"364\u00d77=2548\", \"828\u00d74=3312\"),\n    @(\"399\u00d77=2793\", \"757\u00d75=3785\"),\n    @(\"904\u00d72=1808\", \"943\u00d73=2829\"),\n    @(\"669\u00d76=4014\", \"947\u00d74=3788\"),\n    @(\"128\u00d74=512\", \"465\u00d75=2325\"),\n    @(\"169\u00d77=1183\", \"959\u00d79=8631\"),\n    @(\"143\u00d74=572\", \"693\u00d73=2079\"),\n    @(\"620\u00d77=4340\", \"772\u00d76=4632\"),\n    @(\"619\u00d76=3714\", \"796\u00d78=6368\"),\n    @(\"305\u00d77=2135\", \"324\u00d73=972\"),\n    @(\"321\u00d79=2889\", \"902\u00d76=5412\"),\n    @(\"857\u00d79=7713\", \"871\u00d75=4355\"),\n    @(\"516\u00d72=1032\", \"188\u00d78=1504\"),\n    @(\"522\u00d77=3654\", \"551\u00d75=2755\"),\n    @(\"265\u00d78=2120\", \"560\u00d76=3360\"),\n    @(\"182\u00d76=1092\", \"376\u00d74=1504\"),\n    @(\"182\u00d72=364\", \"545\u00d75=2725\"),\n    @(\"878\u00d79=7902\", \"354\u00d77=2478\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Text = $oldText\n    $rng.Find.Replacement.Text = $newText\n    # 1 = wdFindContinue, 2 = wdReplaceAll\n    $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n\n"}
